# Applies the cryptos-list update described in the commit:
#   "Updated cryptos list on Sat Nov  9 09:59:32 UTC 2024 with GitHub Actions"
#
# Two coin rows were re-ranked (swapped) and every Price (D) / Volume(1h) (E)
# cell in the table was refreshed with newer scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store a number-looking string (e.g. '201.96')
# as text - exactly like typing it into a cell by hand - instead of silently
# auto-converting it to a numeric value. All Price/Volume columns in this sheet
# are plain text.
$quote = [string][char]39

# --- Re-ranking: row 22 (was SuiNetwork) / row 23 (was Polkadot) swap identities ---
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'

# --- Re-ranking: row 26 (was NEARProtocol) / row 27 (was Dai) swap identities ---
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'

# --- Refreshed Price (D) / Volume(1h) (E) text values ---
$ws.Range('D2').Value = '76.587.38'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '3.045.62'
$ws.Range('E3').Value = '  +4.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = $quote + '201.96'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = $quote + '629.75'
$ws.Range('E6').Value = '  +5.39%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('E9').Value = '  +7.32%  '
$ws.Range('D10').Value = '3.043.79'
$ws.Range('E10').Value = '  +4.45%  '
$ws.Range('D11').Value = $quote + '0.440'
$ws.Range('E11').Value = '  +1.69%  '
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').Value = $quote + '5.14'
$ws.Range('E13').Value = '  +5.77%  '
$ws.Range('D14').Value = '3.602.95'
$ws.Range('E14').Value = '  +4.51%  '
$ws.Range('D15').Value = $quote + '29.59'
$ws.Range('E15').Value = '  +6.85%  '
$ws.Range('D16').Value = '76.481.80'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').Value = '3.037.97'
$ws.Range('E18').Value = '  +4.15%  '
$ws.Range('D19').Value = $quote + '13.49'
$ws.Range('E19').Value = '  +4.17%  '
$ws.Range('D20').Value = $quote + '9.08'
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('D21').Value = $quote + '376.65'
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').Value = $quote + '4.37'
$ws.Range('E22').Value = '  +1.86%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').Value = $quote + '2.30'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = $quote + '73.66'
$ws.Range('E24').Value = '  +3.77%  '
$ws.Range('D25').Value = '3.184.03'
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = $quote + '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = $quote + '4.39'
$ws.Range('E27').Value = '  +4.14%  '
$ws.Range('D28').Value = $quote + '9.99'
$ws.Range('E28').Value = '  +3.44%  '
$ws.Range('E29').Value = '  +4.12%  '
$ws.Range('D30').Value = $quote + '0.999'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +8.94%  '
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('D33').Value = $quote + '514.38'
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('E34').Value = '  +8.25%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = $quote + '20.96'
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('D37').Value = $quote + '163.47'
$ws.Range('E37').Value = '  -1.86%  '
$ws.Range('D38').Value = $quote + '0.386'
$ws.Range('E38').Value = '  +11.34%  '
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('E40').Value = '  +1.96%  '
$ws.Range('D41').Value = $quote + '188.48'
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E44').Value = '  +4.73%  '
$ws.Range('D45').Value = $quote + '1.28'
$ws.Range('E45').Value = '  +6.78%  '
$ws.Range('D46').Value = $quote + '42.05'
$ws.Range('E46').Value = '  +4.64%  '
$ws.Range('D47').Value = $quote + '1.68'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('D48').Value = $quote + '0.731'
$ws.Range('E48').Value = '  +11.26%  '
$ws.Range('E49').Value = '  +4.22%  '
$ws.Range('D50').Value = $quote + '0.608'
$ws.Range('E50').Value = '  +6.28%  '
$ws.Range('E51').Value = '  +5.04%  '
